# added ifo gdp component analysis preprocessing
# Fill in the next diagonal of naive QoQ errors matched to ifoCAST,
# one additional forecast-horizon column per row (rows 15-24).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K15").Value = 0.2348700177716323
$ws.Range("J16").Value = 0.2388379152847414
$ws.Range("I17").Value = 0.3744780054549828
$ws.Range("H18").Value = 0.1336718235993181
$ws.Range("G19").Value = 0.08834060834722172
$ws.Range("F20").Value = 0.02147918641116785
$ws.Range("E21").Value = -0.00810701594554874
$ws.Range("D22").Value = -0.02625767267518964
$ws.Range("C23").Value = -0.04428949692388896
$ws.Range("B24").Value = -0.09587373626955231
